# Add a "pre game" row for setting the bet coin amount (BETCOIN / INPUT
# YOUR BET COIN / NHẬP SỐ TIỀN BẠN MUỐN CƯỢC) as a new localisation row
# right after the existing WIN/LOSE rows, then move the selection/scroll
# position to the newly added row like Excel does right after an edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New localisation row (row 31): Code / English / Vietnamese
$ws.Range("A31").Value = "BETCOIN"
$ws.Range("B31").Value = "INPUT YOUR BET COIN: "
$ws.Range("C31").Value = "NHẬP SỐ TIỀN BẠN MUỐN CƯỢC: "

# Scroll the window so row 13 is at the top and select the newly added
# cell, matching where the author ended up after typing the new row.
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C31").Select()
